$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns keep their original text storage (avoid numeric auto-conversion)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "43.076.89"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "2.564.71"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "316.98"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "96.66"
$ws.Range("E6").Value = "  +1.42%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.541"
$ws.Range("E9").Value = "  +2.40%  "
$ws.Range("D10").Value = "35.56"
$ws.Range("E10").Value = "  -1.55%  "
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").Value = "7.45"
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.108"
$ws.Range("E13").Value = "  -4.50%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.958.63"
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("D15").Value = "2.574.22"
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("D16").Value = "15.03"
$ws.Range("E16").Value = "  -2.80%  "
$ws.Range("D17").Value = "0.846"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "43.106.17"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").Value = "6.86"
$ws.Range("E19").Value = "  +4.51%  "
$ws.Range("D20").Value = "12.59"
$ws.Range("E20").Value = "  -3.19%  "
$ws.Range("D21").Value = "0.0₃0962"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "69.52"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("D23").Value = "252.95"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").Value = "2.95"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("E25").Value = "  +2.32%  "
$ws.Range("D26").Value = "26.81"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").Value = "40.01"
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("D30").Value = "10.22"
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").Value = "5.82"
$ws.Range("E31").Value = "  -4.06%  "
$ws.Range("D32").Value = "153.78"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("D33").Value = "3.43"
$ws.Range("E33").Value = "  +4.46%  "
$ws.Range("D34").Value = "2.13"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("D35").Value = "0.0807"
$ws.Range("E35").Value = "  +2.88%  "
$ws.Range("E36").Value = "  +2.82%  "
$ws.Range("D37").Value = "19.02"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("E39").Value = "  +4.24%  "
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").Value = "22.56"
$ws.Range("E41").Value = "  -4.47%  "
$ws.Range("D42").Value = "3.91"
$ws.Range("E42").Value = "  +2.69%  "
$ws.Range("E43").Value = "  +1.19%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "3.27"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("D46").Value = "1.996.12"
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("D47").Value = "9.02"
$ws.Range("E47").Value = "  +2.45%  "
$ws.Range("D48").Value = "83.86"
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("D49").Value = "2.814.08"
$ws.Range("E49").Value = "  +1.38%  "
$ws.Range("D50").Value = "74.21"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "104.15"
$ws.Range("E51").Value = "  +1.53%  "
